$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3009517224809315
$ws.Range("D2").Value = 0.03176124224049381
$ws.Range("E2").Value = 0.1042256271206337
$ws.Range("F2").Value = 1.312523548925626
$ws.Range("G2").Value = 0.002443516547849659
$ws.Range("L2").Value = 0.08654123832794625
$ws.Range("M2").Value = 1.424909482253568
$ws.Range("N2").Value = 1.715880661406516
$ws.Range("O2").Value = 4.681263349677351
$ws.Range("C3").Value = 0.3026857630957807
$ws.Range("D3").Value = 0.03226518995492889
$ws.Range("E3").Value = 0.1062488318439567
$ws.Range("F3").Value = 1.275640392754909
$ws.Range("G3").Value = 0.002448563445509227
$ws.Range("L3").Value = 0.08754708226249974
$ws.Range("M3").Value = 1.293531915880166
$ws.Range("N3").Value = 1.585537508725679
$ws.Range("O3").Value = 4.567228093055689
$ws.Range("C4").Value = 0.3039564342589927
$ws.Range("D4").Value = 0.03259250525759683
$ws.Range("E4").Value = 0.1075626156684795
$ws.Range("F4").Value = 1.253970336737254
$ws.Range("G4").Value = 0.002451826644926911
$ws.Range("L4").Value = 0.0882098158246869
$ws.Range("M4").Value = 1.212821027657299
$ws.Range("N4").Value = 1.505756825956809
$ws.Range("O4").Value = 4.500710647842254
$ws.Range("C5").Value = 0.3045258485037436
$ws.Range("D5").Value = 0.03273037043304505
$ws.Range("E5").Value = 0.1081159277367525
$ws.Range("F5").Value = 1.245383634849162
$ws.Range("G5").Value = 0.002453197902038267
$ws.Range("L5").Value = 0.08849124951990106
$ws.Range("M5").Value = 1.179921701756399
$ws.Range("N5").Value = 1.473312265058979
$ws.Range("O5").Value = 4.474479253443064
$ws.Range("C6").Value = 0.3046235107406545
$ws.Range("D6").Value = 0.03275353302463913
$ws.Range("E6").Value = 0.1082088864720347
$ws.Range("F6").Value = 1.243972514595285
$ws.Range("G6").Value = 0.002453428107490555
$ws.Range("L6").Value = 0.08853866814606448
$ws.Range("M6").Value = 1.174458316800497
$ws.Range("N6").Value = 1.467929010128586
$ws.Range("O6").Value = 4.47017624552592
$ws.Range("C7").Value = 0.3039639048753031
$ws.Range("D7").Value = 0.03259434643324077
$ws.Range("E7").Value = 0.1075700052819895
$ws.Range("F7").Value = 1.253853547413485
$ws.Range("G7").Value = 0.002451844970105421
$ws.Range("L7").Value = 0.08821356530779667
$ws.Range("M7").Value = 1.212377369270399
$ws.Range("N7").Value = 1.505318992175035
$ws.Range("O7").Value = 4.500353346213558
$ws.Range("C8").Value = 0.3015067595107155
$ws.Range("D8").Value = 0.03193128029298897
$ws.Range("E8").Value = 0.1049083515165012
$ws.Range("F8").Value = 1.299602680159026
$ws.Range("G8").Value = 0.002445222685382505
$ws.Range("L8").Value = 0.08687869213064658
$ws.Range("M8").Value = 1.379621020207892
$ws.Range("N8").Value = 1.670888512148537
$ws.Range("O8").Value = 4.641214200906461
$ws.Range("C9").Value = 0.2983311481018234
$ws.Range("D9").Value = 0.03077370177798144
$ws.Range("E9").Value = 0.1002586407527257
$ws.Range("F9").Value = 1.397137411490888
$ws.Range("G9").Value = 0.002433534268498474
$ws.Range("L9").Value = 0.08461860541090616
$ws.Range("M9").Value = 1.707151716691129
$ws.Range("N9").Value = 1.997413350700583
$ws.Range("O9").Value = 4.945477485175161
$ws.Range("C10").Value = 0.297011828438059
$ws.Range("D10").Value = 0.03001123383818971
$ws.Range("E10").Value = 0.09719282262930218
$ws.Range("F10").Value = 1.473676010358815
$ws.Range("G10").Value = 0.002425728968258126
$ws.Range("L10").Value = 0.08317534049809794
$ws.Range("M10").Value = 1.947441446484419
$ws.Range("N10").Value = 2.238268496456499
$ws.Range("O10").Value = 5.186501659324961
$ws.Range("C11").Value = 0.2966344735041275
$ws.Range("D11").Value = 0.02968370105606866
$ws.Range("E11").Value = 0.09587481267580866
$ws.Range("F11").Value = 1.509579998797705
$ws.Range("G11").Value = 0.002422346049547026
$ws.Range("L11").Value = 0.08256577186332947
$ws.Range("M11").Value = 2.056664510312544
$ws.Range("N11").Value = 2.348014863582193
$ws.Range("O11").Value = 5.300032853834864
$ws.Range("C12").Value = 0.2965238381149504
$ws.Range("D12").Value = 0.02956247077993979
$ws.Range("E12").Value = 0.09538679562003249
$ws.Range("F12").Value = 1.523333962904474
$ws.Range("G12").Value = 0.002421089001131804
$ws.Range("L12").Value = 0.0823416885697057
$ws.Range("M12").Value = 2.098010279931231
$ws.Range("N12").Value = 2.389595668641221
$ws.Range("O12").Value = 5.343589814317227
$ws.Range("C13").Value = 0.2965462270790766
$ws.Range("D13").Value = 0.02958845504452512
$ws.Range("E13").Value = 0.09549140481865859
$ws.Range("F13").Value = 1.520364750423852
$ws.Range("G13").Value = 0.002421358664249621
$ws.Range("L13").Value = 0.08238964896779066
$ws.Range("M13").Value = 2.089106415451994
$ws.Range("N13").Value = 2.380639567054686
$ws.Range("O13").Value = 5.334183810566117
$ws.Range("C14").Value = 0.2966247236366826
$ws.Range("D14").Value = 0.02967367114222341
$ws.Range("E14").Value = 0.09583444070153835
$ws.Range("F14").Value = 1.510708371228247
$ws.Range("G14").Value = 0.00242224215116897
$ws.Range("L14").Value = 0.08254720122075909
$ws.Range("M14").Value = 2.060066357884466
$ws.Range("N14").Value = 2.351435316567688
$ws.Range("O14").Value = 5.303604949551982
$ws.Range("C15").Value = 0.2966770129407195
$ws.Range("D15").Value = 0.02972623365438487
$ws.Range("E15").Value = 0.09604600552264619
$ws.Range("F15").Value = 1.504814174853919
$ws.Range("G15").Value = 0.002422786432928008
$ws.Range("L15").Value = 0.08264458495294491
$ws.Range("M15").Value = 2.042276511266351
$ws.Range("N15").Value = 2.333549658016523
$ws.Range("O15").Value = 5.28494829596832
$ws.Range("C16").Value = 0.2970409918124517
$ws.Range("D16").Value = 0.03003302941846542
$ws.Range("E16").Value = 0.09728050537591226
$ws.Range("F16").Value = 1.471351602948999
$ws.Range("G16").Value = 0.002425953414149558
$ws.Range("L16").Value = 0.08321612185823568
$ws.Range("M16").Value = 1.940301537803919
$ws.Range("N16").Value = 2.231099623502757
$ws.Range("O16").Value = 5.179160883265524
$ws.Range("C17").Value = 0.2973215045507658
$ws.Range("D17").Value = 0.03022620151733335
$ws.Range("E17").Value = 0.09805750667384006
$ws.Range("F17").Value = 1.451102809326827
$ws.Range("G17").Value = 0.002427939124906749
$ws.Range("L17").Value = 0.08357876770572048
$ws.Range("M17").Value = 1.877719554253787
$ws.Range("N17").Value = 2.168293358642188
$ws.Range("O17").Value = 5.115263852230669
$ws.Range("C18").Value = 0.2975038066678906
$ws.Range("D18").Value = 0.03033912596727539
$ws.Range("E18").Value = 0.09851162775894062
$ws.Range("F18").Value = 1.43955837325413
$ws.Range("G18").Value = 0.002429097050901593
$ws.Range("L18").Value = 0.08379177445365471
$ws.Range("M18").Value = 1.841716097674691
$ws.Range("N18").Value = 2.132186028664592
$ws.Range("O18").Value = 5.07887756864767
$ws.Range("C19").Value = 0.2975691239988691
$ws.Range("D19").Value = 0.03037767166581951
$ws.Range("E19").Value = 0.09866662184936414
$ws.Range("F19").Value = 1.435667115500863
$ws.Range("G19").Value = 0.002429491821831295
$ws.Range("L19").Value = 0.08386465470322335
$ws.Range("M19").Value = 1.829524649500655
$ws.Range("N19").Value = 2.1199637692539
$ws.Range("O19").Value = 5.066620406709774
$ws.Range("C20").Value = 0.2972894727617899
$ws.Range("D20").Value = 0.03020544981398565
$ws.Range("E20").Value = 0.09797404684582256
$ws.Range("F20").Value = 1.453247745134462
$ws.Range("G20").Value = 0.002427726108705832
$ws.Range("L20").Value = 0.0835397057929228
$ws.Range("M20").Value = 1.884382353119918
$ws.Range("N20").Value = 2.174977442232091
$ws.Range("O20").Value = 5.122027918551339
$ws.Range("C21").Value = 0.2966007899762531
$ws.Range("D21").Value = 0.02964856496435786
$ws.Range("E21").Value = 0.09573338140600107
$ws.Range("F21").Value = 1.513540384705351
$ws.Range("G21").Value = 0.002421981999780416
$ws.Range("L21").Value = 0.08250074125545126
$ws.Range("M21").Value = 2.068596541891282
$ws.Range("N21").Value = 2.360012738277874
$ws.Range("O21").Value = 5.312571315306798
$ws.Range("C22").Value = 0.2963388403708365
$ws.Range("D22").Value = 0.02930093510974352
$ws.Range("E22").Value = 0.0943336168456792
$ws.Range("F22").Value = 1.553866297344626
$ws.Range("G22").Value = 0.002418367661157329
$ws.Range("L22").Value = 0.08186104157371688
$ws.Range("M22").Value = 2.188904829264374
$ws.Range("N22").Value = 2.481071957696315
$ws.Range("O22").Value = 5.440399496426608
$ws.Range("C23").Value = 0.296461358990399
$ws.Range("D23").Value = 0.02948497063503375
$ws.Range("E23").Value = 0.09507476199660481
$ws.Range("F23").Value = 1.532258739864943
$ws.Range("G23").Value = 0.002420283956742058
$ws.Range("L23").Value = 0.08219886592309678
$ws.Range("M23").Value = 2.124702675870793
$ws.Range("N23").Value = 2.41644986278493
$ws.Range("O23").Value = 5.371871488155364
$ws.Range("C24").Value = 0.297303888841185
$ws.Range("D24").Value = 0.03021482584494528
$ws.Range("E24").Value = 0.09801175594991429
$ws.Range("F24").Value = 1.452277718071471
$ws.Range("G24").Value = 0.002427822362529082
$ws.Range("L24").Value = 0.08355735161409328
$ws.Range("M24").Value = 1.881370177083113
$ws.Range("N24").Value = 2.171955565181918
$ws.Range("O24").Value = 5.118968797658567
$ws.Range("C25").Value = 0.299013095977017
$ws.Range("D25").Value = 0.03107147770247387
$ws.Range("E25").Value = 0.1014552064722563
$ws.Range("F25").Value = 1.369902899090064
$ws.Range("G25").Value = 0.002436558285892938
$ws.Range("L25").Value = 0.08519181771271889
$ws.Range("M25").Value = 1.618601066279183
$ws.Range("N25").Value = 1.908899810296077
$ws.Range("O25").Value = 4.860124672401298
